# Generate Report for Handoff
#
# Refresh the localization-status report: the six "Ready for handoff" rows
# (5ec4b178, 66ee08bb, 801b26c0, b22fe8f6, be39082d, e34772d8 -> rows 8-12,14)
# got a new handoff pass, so their handoff timestamps move forward a few
# seconds and they pick up the "ht" (handoff type) priority marker on the
# per-language sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 11, 12, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G
    $wsOverview.Range("G$r").Value = "2016-08-29 20:32:30"

    # zh-cn sheet: Priority column E, "Latest Handoff Datetime" column H
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-29 20:32:25"

    # de-de sheet: Priority column E, "Latest Handoff Datetime" column H
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-29 20:32:30"
}
